$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A29").Value = 111277538
$ws.Range("Q29").Value = 569784.6764437903
$ws.Range("R29").Value = 6992856.400962653
$ws.Range("A30").Value = 111279094
$ws.Range("B30").Value = 89416
$ws.Range("D30").Value = "LC"
$ws.Range("E30").Value = 1205
$ws.Range("F30").Value = "Stor aspticka"
$ws.Range("G30").Value = "Phellinus populicola"
$ws.Range("H30").Value = "Niemelä"
$ws.Range("Q30").Value = 569279.6199819668
$ws.Range("R30").Value = 6992811.114809629
$ws.Range("Z30").Value = "20:26"
$ws.Range("AB30").Value = "20:26"
$ws.Range("A32").Value = 111277392
$ws.Range("B32").Value = 89845
$ws.Range("E32").Value = 1209
$ws.Range("F32").Value = "Rynkskinn"
$ws.Range("G32").Value = "Phlebia centrifuga"
$ws.Range("H32").Value = "P.Karst."
$ws.Range("Q32").Value = 569750.3053765292
$ws.Range("R32").Value = 6992912.817455334
$ws.Range("Z32").Value = "00:00"
$ws.Range("AB32").Value = "00:00"
$ws.Range("A33").Value = 111278872
$ws.Range("B33").Value = 78578
$ws.Range("E33").Value = 6458
$ws.Range("F33").Value = "Lunglav"
$ws.Range("G33").Value = "Lobaria pulmonaria"
$ws.Range("H33").Value = "(L.) Hoffm."
$ws.Range("P33").Value = "Ragunda, Jmt"
$ws.Range("Q33").Value = 569296.7869269754
$ws.Range("R33").Value = 6992794.243538878
$ws.Range("A34").Value = 111278217
$ws.Range("B34").Value = 89686
$ws.Range("D34").Value = "NT"
$ws.Range("E34").Value = 658
$ws.Range("F34").Value = "Rosenticka"
$ws.Range("G34").Value = "Rhodofomes roseus"
$ws.Range("H34").Value = "(Alb. & Schwein.) Kotl. & Pouzar"
$ws.Range("P34").Value = "Singsån, Jmt"
$ws.Range("Q34").Value = 569671.7019483433
$ws.Range("R34").Value = 6993040.858867787
$ws.Range("A35").Value = 111277448
$ws.Range("Q35").Value = 569750.3053765292
$ws.Range("R35").Value = 6992912.817455334
$ws.Range("A36").Value = 111279409
$ws.Range("B36").Value = 96348
$ws.Range("D36").Value = "VU"
$ws.Range("E36").Value = 220787
$ws.Range("F36").Value = "Knärot"
$ws.Range("G36").Value = "Goodyera repens"
$ws.Range("H36").Value = "(L.) R. Br."
$ws.Range("Q36").Value = 569443.239979364
$ws.Range("R36").Value = 6992913.042043422
$ws.Range("Z36").Value = "21:39"
$ws.Range("AB36").Value = "21:39"
$ws.Range("A37").Value = 111278492
$ws.Range("B37").Value = 89686
$ws.Range("E37").Value = 658
$ws.Range("F37").Value = "Rosenticka"
$ws.Range("G37").Value = "Rhodofomes roseus"
$ws.Range("H37").Value = "(Alb. & Schwein.) Kotl. & Pouzar"
$ws.Range("Q37").Value = 569641.4769454591
$ws.Range("R37").Value = 6992967.635971196
$ws.Range("A38").Value = 111277903
$ws.Range("B38").Value = 96348
$ws.Range("D38").Value = "VU"
$ws.Range("E38").Value = 220787
$ws.Range("F38").Value = "Knärot"
$ws.Range("G38").Value = "Goodyera repens"
$ws.Range("H38").Value = "(L.) R. Br."
$ws.Range("P38").Value = "Präststranden, Jmt"
$ws.Range("Q38").Value = 569897.0842333297
$ws.Range("R38").Value = 6993078.813114846
$ws.Range("Z38").Value = "20:24"
$ws.Range("AB38").Value = "20:24"
$ws.Range("A39").Value = 111277950
$ws.Range("Q39").Value = 569874.8142812594
$ws.Range("R39").Value = 6993100.559414167
$ws.Range("Z39").Value = "20:26"
$ws.Range("AB39").Value = "20:26"
$ws.Range("A40").Value = 111277552
$ws.Range("Q40").Value = 569770.841244747
$ws.Range("R40").Value = 6992866.083226931
$ws.Range("Z40").Value = "00:00"
$ws.Range("AB40").Value = "00:00"
$ws.Range("A41").Value = 111277633
$ws.Range("Q41").Value = 569802.0407188418
$ws.Range("R41").Value = 6992830.464391444
